$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'243.56"
$ws.Range("G2").Formula = "'11"
$ws.Range("D3").Formula = "'23.86"
$ws.Range("G3").Formula = "'11"
$ws.Range("D4").Formula = "'5.262"
$ws.Range("G4").Formula = "'11"
$ws.Range("D5").Formula = "'0.05816"
$ws.Range("G5").Formula = "'11"
$ws.Range("D6").Formula = "'6.472"
$ws.Range("G6").Formula = "'11"
$ws.Range("D7").Formula = "'3.329"
$ws.Range("G7").Formula = "'11"
$ws.Range("D8").Formula = "'0.8075"
$ws.Range("G8").Formula = "'11"
$ws.Range("D9").Formula = "'0.8754"
$ws.Range("G9").Formula = "'11"
$ws.Range("D10").Formula = "'0.1384"
$ws.Range("G10").Formula = "'11"
$ws.Range("D11").Formula = "'0.07263"
$ws.Range("G11").Formula = "'11"
$ws.Range("D12").Formula = "'0.03082"
$ws.Range("G12").Formula = "'11"
$ws.Range("D13").Formula = "'0.03054"
$ws.Range("G13").Formula = "'11"
$ws.Range("D14").Formula = "'0.09319"
$ws.Range("G14").Formula = "'11"
$ws.Range("D15").Formula = "'3.861"
$ws.Range("G15").Formula = "'11"
$ws.Range("D16").Formula = "'0.001552"
$ws.Range("G16").Formula = "'11"
$ws.Range("D17").Formula = "'0.04702"
$ws.Range("G17").Formula = "'11"
$ws.Range("D18").Formula = "'0.0006028"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Formula = "'11"
$ws.Range("D19").Formula = "'0.006152"
$ws.Range("G19").Formula = "'11"
$ws.Range("D20").Formula = "'0.001264"
$ws.Range("G20").Formula = "'11"
$ws.Range("D21").Formula = "'0.004595"
$ws.Range("G21").Formula = "'11"
$ws.Range("D22").Formula = "'0.00008695"
$ws.Range("G22").Formula = "'11"
$ws.Range("D23").Formula = "'3.561"
$ws.Range("G23").Formula = "'11"
$ws.Range("D24").Formula = "'2.184"
$ws.Range("G24").Formula = "'11"
$ws.Range("D25").Formula = "'0.3197"
$ws.Range("G25").Formula = "'11"
$ws.Range("G26").Formula = "'11"
$ws.Range("G27").Formula = "'11"
$ws.Range("D28").Formula = "'0.0002343"
$ws.Range("G28").Formula = "'11"
$ws.Range("G29").Formula = "'11"
$ws.Range("G30").Formula = "'11"
$ws.Range("G31").Formula = "'11"
$ws.Range("G32").Formula = "'11"
$ws.Range("G33").Formula = "'11"
$ws.Range("G34").Formula = "'11"
$ws.Range("G35").Formula = "'11"
$ws.Range("G36").Formula = "'11"
$ws.Range("G37").Formula = "'11"
$ws.Range("G38").Formula = "'11"
$ws.Range("G39").Formula = "'11"
$ws.Range("D40").Formula = "'0.03782"
$ws.Range("G40").Formula = "'11"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Formula = "'0.1055"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").Formula = "'11"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Formula = "'0.002608"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Formula = "'11"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Formula = "'0.003243"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("G43").Formula = "'11"
$ws.Range("G44").Formula = "'11"
$ws.Range("D45").Formula = "'0.00005531"
$ws.Range("G45").Formula = "'11"
$ws.Range("G46").Formula = "'11"
$ws.Range("D47").Formula = "'0.5499"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Formula = "'11"
$ws.Range("D48").Formula = "'0.01415"
$ws.Range("G48").Formula = "'11"
$ws.Range("D49").Formula = "'0.00002099"
$ws.Range("G49").Formula = "'11"
$ws.Range("D50").Formula = "'0.0001999"
$ws.Range("G50").Formula = "'11"
$ws.Range("G51").Formula = "'11"
